# MazeDimensions_human.xlsx - add dynamic erosion radius, depending on solver and size
# Update the dimension table values that shifted with the new "human" erosion-radius
# calibration, widen the data columns so the longer numbers are readable, and leave the
# cursor parked on the next empty row (F6) like the author's last saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated measurement pairs (column C/D/E/F, rows 2-4) ---
$ws.Range("D2").Value = "3.31, 1.3"
$ws.Range("C2").Value = "6.34, 2.15"
$ws.Range("D3").Value = "3.24, 1.3"
$ws.Range("E3").Value = "3.24, 2.00"
$ws.Range("E2").Value = "3.31, 1.93"
$ws.Range("F3").Value = "2.18, 0.91"
$ws.Range("F4").Value = "4.36, 2.19"
$ws.Range("D4").Value = "6.52, 2.63"
$ws.Range("E4").Value = "6.52, 3.76"
$ws.Range("F2").Value = "2.22, 0.91"

# --- Widen the label / value columns to fit the new content ---
# (input values pre-compensate for the host's internal pixel-grid rounding of
# ColumnWidth so the saved OOXML <col width="..."> lands as close as possible
# to the target 31.33203125 / 24.33203125 / 22.33203125 / 19.109375 chars)
$ws.Columns.Item(2).ColumnWidth = 30.41796875
$ws.Columns.Item(3).ColumnWidth = 23.41796875
$ws.Columns.Item(4).ColumnWidth = 21.41796875
$ws.Columns.Item(5).ColumnWidth = 18.251302083333332

# --- Leave selection on the next free row, matching the saved view state ---
$ws.Range("F6").Select()
